$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bulk-update column B (epoch accuracy values) for rows 2-118
$bVals = New-Object 'object[,]' 117,1
$bVals[0,0] = 0.921875
$bVals[1,0] = 0.84375
$bVals[2,0] = 0.703125
$bVals[3,0] = 0.703125
$bVals[4,0] = 0.65625
$bVals[5,0] = 0.515625
$bVals[6,0] = 0.484375
$bVals[7,0] = 0.515625
$bVals[8,0] = 0.4375
$bVals[9,0] = 0.421875
$bVals[10,0] = 0.5
$bVals[11,0] = 0.375
$bVals[12,0] = 0.5625
$bVals[13,0] = 0.484375
$bVals[14,0] = 0.421875
$bVals[15,0] = 0.40625
$bVals[16,0] = 0.359375
$bVals[17,0] = 0.359375
$bVals[18,0] = 0.546875
$bVals[19,0] = 0.359375
$bVals[20,0] = 0.359375
$bVals[21,0] = 0.34375
$bVals[22,0] = 0.296875
$bVals[23,0] = 0.3125
$bVals[24,0] = 0.296875
$bVals[25,0] = 0.28125
$bVals[26,0] = 0.296875
$bVals[27,0] = 0.296875
$bVals[28,0] = 0.4375
$bVals[29,0] = 0.3125
$bVals[30,0] = 0.3125
$bVals[31,0] = 0.359375
$bVals[32,0] = 0.296875
$bVals[33,0] = 0.28125
$bVals[34,0] = 0.265625
$bVals[35,0] = 0.296875
$bVals[36,0] = 0.296875
$bVals[37,0] = 0.3125
$bVals[38,0] = 0.328125
$bVals[39,0] = 0.328125
$bVals[40,0] = 0.359375
$bVals[41,0] = 0.359375
$bVals[42,0] = 0.296875
$bVals[43,0] = 0.28125
$bVals[44,0] = 0.265625
$bVals[45,0] = 0.265625
$bVals[46,0] = 0.25
$bVals[47,0] = 0.21875
$bVals[48,0] = 0.21875
$bVals[49,0] = 0.203125
$bVals[50,0] = 0.21875
$bVals[51,0] = 0.21875
$bVals[52,0] = 0.21875
$bVals[53,0] = 0.21875
$bVals[54,0] = 0.21875
$bVals[55,0] = 0.21875
$bVals[56,0] = 0.21875
$bVals[57,0] = 0.21875
$bVals[58,0] = 0.21875
$bVals[59,0] = 0.21875
$bVals[60,0] = 0.21875
$bVals[61,0] = 0.21875
$bVals[62,0] = 0.21875
$bVals[63,0] = 0.21875
$bVals[64,0] = 0.203125
$bVals[65,0] = 0.203125
$bVals[66,0] = 0.203125
$bVals[67,0] = 0.203125
$bVals[68,0] = 0.203125
$bVals[69,0] = 0.21875
$bVals[70,0] = 0.21875
$bVals[71,0] = 0.21875
$bVals[72,0] = 0.1875
$bVals[73,0] = 0.34375
$bVals[74,0] = 0.203125
$bVals[75,0] = 0.21875
$bVals[76,0] = 0.1875
$bVals[77,0] = 0.234375
$bVals[78,0] = 0.21875
$bVals[79,0] = 0.234375
$bVals[80,0] = 0.21875
$bVals[81,0] = 0.203125
$bVals[82,0] = 0.203125
$bVals[83,0] = 0.203125
$bVals[84,0] = 0.203125
$bVals[85,0] = 0.203125
$bVals[86,0] = 0.203125
$bVals[87,0] = 0.203125
$bVals[88,0] = 0.21875
$bVals[89,0] = 0.21875
$bVals[90,0] = 0.21875
$bVals[91,0] = 0.21875
$bVals[92,0] = 0.203125
$bVals[93,0] = 0.203125
$bVals[94,0] = 0.203125
$bVals[95,0] = 0.203125
$bVals[96,0] = 0.203125
$bVals[97,0] = 0.203125
$bVals[98,0] = 0.203125
$bVals[99,0] = 0.203125
$bVals[100,0] = 0.203125
$bVals[101,0] = 0.328125
$bVals[102,0] = 0.3125
$bVals[103,0] = 0.171875
$bVals[104,0] = 0.265625
$bVals[105,0] = 0.3125
$bVals[106,0] = 0.25
$bVals[107,0] = 0.234375
$bVals[108,0] = 0.203125
$bVals[109,0] = 0.1875
$bVals[110,0] = 0.25
$bVals[111,0] = 0.1875
$bVals[112,0] = 0.28125
$bVals[113,0] = 0.21875
$bVals[114,0] = 0.265625
$bVals[115,0] = 0.265625
$bVals[116,0] = 0.180327868852459
$ws.Range("B2:B118").Value = $bVals

# Refresh the repr() placeholder text in column A for the re-run rows (102-118),
# which now carries a new object id, and append the 3 new epoch rows (116-118)
$aVals = New-Object 'object[,]' 17,1
$aVals[0,0] = "<__main__.DisplayOutputs object at 0x7f7eac2c99d0>"
$aVals[1,0] = "<__main__.DisplayOutputs object at 0x7f7eac2c99d0>"
$aVals[2,0] = "<__main__.DisplayOutputs object at 0x7f7eac2c99d0>"
$aVals[3,0] = "<__main__.DisplayOutputs object at 0x7f7eac2c99d0>"
$aVals[4,0] = "<__main__.DisplayOutputs object at 0x7f7eac2c99d0>"
$aVals[5,0] = "<__main__.DisplayOutputs object at 0x7f7eac2c99d0>"
$aVals[6,0] = "<__main__.DisplayOutputs object at 0x7f7eac2c99d0>"
$aVals[7,0] = "<__main__.DisplayOutputs object at 0x7f7eac2c99d0>"
$aVals[8,0] = "<__main__.DisplayOutputs object at 0x7f7eac2c99d0>"
$aVals[9,0] = "<__main__.DisplayOutputs object at 0x7f7eac2c99d0>"
$aVals[10,0] = "<__main__.DisplayOutputs object at 0x7f7eac2c99d0>"
$aVals[11,0] = "<__main__.DisplayOutputs object at 0x7f7eac2c99d0>"
$aVals[12,0] = "<__main__.DisplayOutputs object at 0x7f7eac2c99d0>"
$aVals[13,0] = "<__main__.DisplayOutputs object at 0x7f7eac2c99d0>"
$aVals[14,0] = "<__main__.DisplayOutputs object at 0x7f7eac2c99d0>"
$aVals[15,0] = "<__main__.DisplayOutputs object at 0x7f7eac2c99d0>"
$aVals[16,0] = "<__main__.DisplayOutputs object at 0x7f7eac2c99d0>"
$ws.Range("A102:A118").Value = $aVals
